# Auto-generated edit script applying the cryptos.xlsx diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.152.21'
$ws.Range("E2").Value = '  -4.30%  '
$ws.Range("D3").Value = '2.966.64'
$ws.Range("E3").Value = '  -6.40%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '123.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.74%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '2.965.30'
$ws.Range("E8").Value = '  -6.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.43%  '
$ws.Range("E10").Value = '  -6.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.10'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.434'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.37%  '
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("D16").Value = '3.466.99'
$ws.Range("E16").Value = '  -5.92%  '
$ws.Range("D17").Value = '60.198.97'
$ws.Range("E17").Value = '  -4.20%  '
$ws.Range("D18").Value = '2.976.23'
$ws.Range("E18").Value = '  -5.92%  '
$ws.Range("E19").Value = '  -6.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '423.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.657'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.56%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.11'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.61%  '
$ws.Range("E32").Value = '  -11.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0918'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.23'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.59%  '
$ws.Range("E35").Value = '  -8.91%  '
$ws.Range("E36").Value = '  -4.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '49.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.14%  '
$ws.Range("D38").Value = '0.0₃0646'
$ws.Range("E38").Value = '  -8.20%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.24%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0354'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.55%  '
$ws.Range("E41").Value = '  -3.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '376.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.24%  '
$ws.Range("D43").Value = '2.627.76'
$ws.Range("E43").Value = '  -5.92%  '
$ws.Range("E44").Value = '  -9.23%  '
$ws.Range("E46").Value = '  -7.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.31'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.56%  '
$ws.Range("E48").Value = '  -7.78%  '
$ws.Range("E49").Value = '  -4.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.97%  '
